$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '59.980.79'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.389.45'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.24%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '558.94'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '134.03'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.43%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.586'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.61%  '
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.62'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.37%  '
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.344'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.92%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '24.46'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -4.20%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.818.29'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.06%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '59.896.88'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('E16').Value = '  -0.59%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.392.32'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.03%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.10'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.19%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.49'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.85%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '322.02'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.15%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.74'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.17%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '64.27'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.58%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.173'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('B25').Value = 'Binance-PegBSC-USD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.01'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.11%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.47'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.13%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.38'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.66%  '
$ws.Range('E28').Value = '  +1.59%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0₃0764'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.59%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '170.62'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.08'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.08'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +6.60%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.399'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.40%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '18.20'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.43%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.32'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.47%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.15'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.96%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.59'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.15%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '319.45'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.50%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '38.69'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -2.34%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '146.92'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +5.90%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.53'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.83%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0967'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '19.77'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.88%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0512'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.36%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.572'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.36%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0218'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.73%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '11.07'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.54'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.23%  '
$ws.Range('E51').Value = '  +0.28%  '
